$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of rows 3 to 6 in column A (remove the now-unused shared strings)
$ws.Range("A3:A6").ClearContents()

# Update the selection to C22
$ws.Range("C22").Select()
